# Update Logic component class diagram
# Rename AddressBookParser to BattleshipParser.
#
# The "AddressBookParser" class box on the diagram is a two-line shape
# (line 1: "AddressBook", line 2: "Parser"). Only the first line needs to
# change to "Battleship" so the box now reads "BattleshipParser".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$renamed = $false

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tf = $shp.TextFrame
    if (-not $tf.HasText) { continue }

    $tr = $tf.TextRange
    if ($tr.Text -eq "AddressBook" + [char]13 + "Parser") {
        # Replace just the "AddressBook" portion (first paragraph, minus the
        # trailing paragraph mark) with "Battleship", keeping the "Parser"
        # line on the second paragraph untouched so the box reads
        # "BattleshipParser".
        $firstPara = $tr.Paragraphs(1, 1)
        $word = $tr.Characters(1, $firstPara.Length - 1)
        $word.Text = "Battleship"
        $renamed = $true
        break
    }
}

if (-not $renamed) {
    throw "Could not find the AddressBookParser shape to rename"
}
